# Generate Report for Handoff
#
# The localization job for "a4b0715a-fcaf-48bf-939d-1a3b8497fb0e.md" has
# finished translation and is ready to be handed off. Update its status
# row on every sheet (the per-locale "zh-cn"/"de-de" detail sheets and the
# "Overview" summary sheet) to reflect the new "Ready for handoff" state,
# the new machine-translation ("mt") priority, and the refreshed handoff
# timestamps.

$wb = $excel.ActiveWorkbook

$statusOld  = "In Translation"
$statusNew  = "Ready for handoff"
$priorityOld = "ht"
$priorityNew = "mt"

# ---------------------------------------------------------------------
# "zh-cn" detail sheet - row 3 is the a4b0715a...md file
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusNew
$wsZh.Range("E3").Value = $priorityNew
$wsZh.Range("H3").Value = "2016-09-03 20:14:36"

# ---------------------------------------------------------------------
# "de-de" detail sheet - row 3 is the a4b0715a...md file
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusNew
$wsDe.Range("E3").Value = $priorityNew
$wsDe.Range("H3").Value = "2016-09-03 20:14:41"

# ---------------------------------------------------------------------
# "Overview" sheet - row 3 is the a4b0715a...md file; E=zh-cn, F=de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = "2016-09-03 20:14:41"

# ---------------------------------------------------------------------
# The wider "Ready for handoff" text (vs. "In Translation") causes the
# Status columns to auto-fit wider in the source workbook. Mirror that by
# widening the corresponding columns: Overview!E:F and the "Status"
# column (C) on each detail sheet.
# ---------------------------------------------------------------------
$wsOverview.Range("E1:F1").ColumnWidth = 16.3
$wsZh.Range("C1").ColumnWidth = 16.3
$wsDe.Range("C1").ColumnWidth = 16.3
